$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2. This pushes the existing row 2 (LUCAS CORTES...)
# down to row 3, and the existing row 3 (the empty styled row) down to row 4.
$ws.Rows.Item(2).Insert()
$ws.Range("A2:E2").Style = "Normal"

# Fill in the new row 2 with the new record's data.
$ws.Cells.Item(2, 1).Value = "LUCAS"

# DEPTO. is a text column (sibling value "221" in row 3 is text too), but a
# plain numeric-looking string like "1" would be auto-coerced to a number.
# Force text storage, then drop the now-unneeded explicit number format so
# the cell keeps the sheet's default (unstyled) look.
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = "1"
$ws.Cells.Item(2, 2).ClearFormats()

$ws.Cells.Item(2, 3).Value = 10
$ws.Cells.Item(2, 4).Value = 50000
$ws.Cells.Item(2, 5).Value = "2021/01/07, 12:31:19"
